$d = $word.ActiveDocument

# Locate the paragraph that ends the first bonus question
# ("Use your data to determine whether the mean or the median
# summarizes the data more meaningfully.") via Find, then insert a
# new sub-bullet answer paragraph right after it.
$findRange = $d.Content
$findRange.Find.Execute("Use your data to determine whether the mean or the median summarizes the data more meaningfully.")
$questionPara = $findRange.Paragraphs(1)

$questionPara.Range.InsertParagraphAfter()

# Re-fetch the freshly inserted paragraph, give it the same list
# style as the other sub-bullet answers (ilvl = 1, numId = 2), and
# fill in its text.
$newPara = $questionPara.Next()
$newPara.Range.ListFormat.ListLevelNumber = 2
$newPara.Range.Text = "For the successful state date, the median summarizes the data more meaningfully due the large outliers that skews the mean. For the failed state data, either the mean or median does a good job summarizing the data as there are very few outliers to skew the mean. "
